$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_2_9_24"
$ws.Range("B2").Value = 0.9973525866107236
$ws.Range("C2").Value = 0.8122600222208161
$ws.Range("D2").Value = 0.9691608515865253
$ws.Range("E2").Value = 0.9970455838252948
$ws.Range("F2").Value = 0.9940532337571081
$ws.Range("G2").Value = 0.01770326552430962
$ws.Range("H2").Value = 1.255418095872513
$ws.Range("I2").Value = 0.05240293613909332
$ws.Range("J2").Value = 0.04093993027686369
$ws.Range("K2").Value = 0.04667144626999713
$ws.Range("L2").Value = 0.1657040364817196
$ws.Range("M2").Value = 0.1330536189823848
$ws.Range("N2").Value = 1.001411953807614
$ws.Range("O2").Value = 0.1387179870014262
$ws.Range("P2").Value = 146.0680123270405
$ws.Range("Q2").Value = 230.1704442429463
$ws.Range("A3").Value = "model_2_9_23"
$ws.Range("B3").Value = 0.9973826592702316
$ws.Range("C3").Value = 0.812163785848439
$ws.Range("D3").Value = 0.9703160096145047
$ws.Range("E3").Value = 0.9971589745632887
$ws.Range("F3").Value = 0.994278392736665
$ws.Range("G3").Value = 0.01750216951170702
$ws.Range("H3").Value = 1.256061628937727
$ws.Range("I3").Value = 0.05044005209446401
$ws.Range("J3").Value = 0.03936865235493263
$ws.Range("K3").Value = 0.04490435222469831
$ws.Range("L3").Value = 0.1727771143498103
$ws.Range("M3").Value = 0.1322957652825933
$ws.Range("N3").Value = 1.001395915055876
$ws.Range("O3").Value = 0.1379278699006612
$ws.Range("P3").Value = 146.0908608672781
$ws.Range("Q3").Value = 230.1932927831839
$ws.Range("A4").Value = "model_2_9_22"
$ws.Range("B4").Value = 0.9974086494768584
$ws.Range("C4").Value = 0.8120559255979419
$ws.Range("D4").Value = 0.9715138207837054
$ws.Range("E4").Value = 0.9972764425320984
$ws.Range("F4").Value = 0.9945117670044527
$ws.Range("G4").Value = 0.01732837288032005
$ws.Range("H4").Value = 1.256782890929451
$ws.Range("I4").Value = 0.04840469037290322
$ws.Range("J4").Value = 0.03774087543778502
$ws.Range("K4").Value = 0.04307278290534412
$ws.Range("L4").Value = 0.1803022654682103
$ws.Range("M4").Value = 0.1316372777002018
$ws.Range("N4").Value = 1.001382053612342
$ws.Range("O4").Value = 0.1372413491386301
$ws.Range("P4").Value = 146.1108201400153
$ws.Range("Q4").Value = 230.2132520559211
$ws.Range("B5").Value = 0.9974295149275144
$ws.Range("C5").Value = 0.8119350801102593
$ws.Range("D5").Value = 0.9727530972059604
$ws.Range("E5").Value = 0.9973978401491324
$ws.Range("F5").Value = 0.9947530993387945
$ws.Range("G5").Value = 0.01718884551570654
$ws.Range("H5").Value = 1.257590985262027
$ws.Range("I5").Value = 0.04629886947462784
$ws.Range("J5").Value = 0.03605864460663114
$ws.Range("K5").Value = 0.0411787570406295
$ws.Range("L5").Value = 0.1883176757760545
$ws.Range("M5").Value = 0.1311062375163994
$ws.Range("N5").Value = 1.001370925371992
$ws.Range("O5").Value = 0.1366877014747985
$ws.Range("P5").Value = 146.1269892441385
$ws.Range("Q5").Value = 230.2294211600443
$ws.Range("A6").Value = "model_2_9_20"
$ws.Range("B6").Value = 0.9974441232596961
$ws.Range("C6").Value = 0.8117996820807794
$ws.Range("D6").Value = 0.9740325767062273
$ws.Range("E6").Value = 0.9975230709968401
$ws.Range("F6").Value = 0.9950021679814866
$ws.Range("G6").Value = 0.01709115953114242
$ws.Range("H6").Value = 1.25849639250861
$ws.Range("I6").Value = 0.04412473413065485
$ws.Range("J6").Value = 0.0343232959385715
$ws.Range("K6").Value = 0.03922401503461317
$ws.Range("L6").Value = 0.196844370190955
$ws.Range("M6").Value = 0.1307331615587354
$ws.Range("N6").Value = 1.001363134261495
$ws.Range("O6").Value = 0.136298742901243
$ws.Range("P6").Value = 146.13838787128
$ws.Range("Q6").Value = 230.2408197871859
$ws.Range("A7").Value = "model_2_9_19"
$ws.Range("B7").Value = 0.997451060348341
$ws.Range("C7").Value = 0.8116479107177935
$ws.Range("D7").Value = 0.9753501810958675
$ws.Range("E7").Value = 0.9976518619625084
$ws.Range("F7").Value = 0.995258513797711
$ws.Range("G7").Value = 0.01704477118743183
$ws.Range("H7").Value = 1.259511288311739
$ws.Range("I7").Value = 0.0418858156702238
$ws.Range("J7").Value = 0.03253861401058391
$ws.Range("K7").Value = 0.03721216027190746
$ws.Range("L7").Value = 0.2059262776550833
$ws.Range("M7").Value = 0.1305556248785621
$ws.Range("N7").Value = 1.001359434480885
$ws.Range("O7").Value = 0.1361136481170432
$ws.Range("P7").Value = 146.143823595022
$ws.Range("Q7").Value = 230.2462555109279
$ws.Range("A8").Value = "model_2_9_18"
$ws.Range("B8").Value = 0.9974488338758669
$ws.Range("C8").Value = 0.8114778150902741
$ws.Range("D8").Value = 0.9767038301346265
$ws.Range("E8").Value = 0.9977840235289528
$ws.Range("F8").Value = 0.9955217271524078
$ws.Range("G8").Value = 0.01705965961911797
$ws.Range("H8").Value = 1.26064871855618
$ws.Range("I8").Value = 0.03958564890874994
$ws.Range("J8").Value = 0.03070722499984101
$ws.Range("K8").Value = 0.03514640765283298
$ws.Range("L8").Value = 0.2155957010449485
$ws.Range("M8").Value = 0.1306126319278422
$ws.Range("N8").Value = 1.001360621932871
$ws.Range("O8").Value = 0.1361730820744319
$ws.Range("P8").Value = 146.1420773782074
$ws.Range("Q8").Value = 230.2445092941132
$ws.Range("A9").Value = "model_2_9_17"
$ws.Range("B9").Value = 0.9974356704061602
$ws.Range("C9").Value = 0.8112871962857764
$ws.Range("D9").Value = 0.9780902773129534
$ws.Range("E9").Value = 0.9979192206702533
$ws.Range("F9").Value = 0.9957911660670781
$ws.Range("G9").Value = 0.01714768380165893
$ws.Range("H9").Value = 1.261923387379575
$ws.Range("I9").Value = 0.03722975042634089
$ws.Range("J9").Value = 0.02883377142689539
$ws.Range("K9").Value = 0.03303179555687073
$ws.Range("L9").Value = 0.2258872171572762
$ws.Range("M9").Value = 0.1309491649521253
$ws.Range("N9").Value = 1.001367642450048
$ws.Range("O9").Value = 0.1365239420062782
$ws.Range("P9").Value = 146.1317843395757
$ws.Range("Q9").Value = 230.2342162554815
$ws.Range("B10").Value = 0.9974095242374942
$ws.Range("C10").Value = 0.811073628491664
$ws.Range("D10").Value = 0.9795056044542044
$ws.Range("E10").Value = 0.998057069818036
$ws.Range("F10").Value = 0.9960660852357549
$ws.Range("G10").Value = 0.01732252335192032
$ws.Range("H10").Value = 1.263351516202194
$ws.Range("I10").Value = 0.03482477812280937
$ws.Range("J10").Value = 0.02692356847469436
$ws.Range("K10").Value = 0.03087417329875187
$ws.Range("L10").Value = 0.2368457797421039
$ws.Range("M10").Value = 0.1316150574665388
$ws.Range("N10").Value = 1.001381587073336
$ws.Range("O10").Value = 0.1372181829436178
$ws.Range("P10").Value = 146.1114953927899
$ws.Range("Q10").Value = 230.2139273086957
$ws.Range("B11").Value = 0.9973680921085936
$ws.Range("C11").Value = 0.8108342647301608
$ws.Range("D11").Value = 0.9809451572764264
$ws.Range("E11").Value = 0.9981971542697761
$ws.Range("F11").Value = 0.996345600086373
$ws.Range("G11").Value = 0.01759958018865615
$ws.Range("H11").Value = 1.264952142777534
$ws.Range("I11").Value = 0.03237864071329558
$ws.Range("J11").Value = 0.02498239047268822
$ws.Range("K11").Value = 0.0286804831822319
$ws.Range("L11").Value = 0.2485104399525094
$ws.Range("M11").Value = 0.1326634093812463
$ws.Range("N11").Value = 1.00140368420875
$ws.Range("O11").Value = 0.138311165369722
$ws.Range("P11").Value = 146.0797604602796
$ws.Range("Q11").Value = 230.1821923761854
$ws.Range("B12").Value = 0.9973087335462102
$ws.Range("C12").Value = 0.8105661355773596
$ws.Range("D12").Value = 0.9824038362714499
$ws.Range("E12").Value = 0.9983389136095475
$ws.Range("F12").Value = 0.9966286558391892
$ws.Range("G12").Value = 0.01799651116863536
$ws.Range("H12").Value = 1.266745123656937
$ws.Range("I12").Value = 0.02990000345655944
$ws.Range("J12").Value = 0.02301800321539228
$ws.Range("K12").Value = 0.02645900333597585
$ws.Range("L12").Value = 0.2609377852513163
$ws.Range("M12").Value = 0.1341510759130741
$ws.Range("N12").Value = 1.001435342108688
$ws.Range("O12").Value = 0.139862164945704
$ws.Range("P12").Value = 146.0351547276739
$ws.Range("Q12").Value = 230.1375866435798
$ws.Range("B13").Value = 0.9972284737104338
$ws.Range("C13").Value = 0.8102656899024906
$ws.Range("D13").Value = 0.9838753359983785
$ws.Range("E13").Value = 0.9984817115265839
$ws.Range("F13").Value = 0.9969140205285284
$ws.Range("G13").Value = 0.01853320906003456
$ws.Range("H13").Value = 1.268754205268211
$ws.Range("I13").Value = 0.02739958077351159
$ws.Range("J13").Value = 0.02103922418716847
$ws.Range("K13").Value = 0.02421940248034003
$ws.Range("L13").Value = 0.2741642074396839
$ws.Range("M13").Value = 0.136136729283594
$ws.Range("N13").Value = 1.001478147354435
$ws.Range("O13").Value = 0.1419323516910761
$ws.Range("P13").Value = 145.97638214362
$ws.Range("Q13").Value = 230.0788140595259
$ws.Range("B14").Value = 0.9971239020965844
$ws.Range("C14").Value = 0.8099291986012556
$ws.Range("D14").Value = 0.9853519110187604
$ws.Range("E14").Value = 0.9986248505401021
$ws.Range("F14").Value = 0.9972002327060814
$ws.Range("G14").Value = 0.0192324799233532
$ws.Range("H14").Value = 1.27100432414897
$ws.Range("I14").Value = 0.02489053397817775
$ws.Range("J14").Value = 0.01905571851741831
$ws.Range("K14").Value = 0.02197315036265383
$ws.Range("L14").Value = 0.2882592708948553
$ws.Range("M14").Value = 0.1386812169089715
$ws.Range("N14").Value = 1.001533918881822
$ws.Range("O14").Value = 0.1445851634224814
$ws.Range("P14").Value = 145.9023095330468
$ws.Range("Q14").Value = 230.0047414489527
$ws.Range("B15").Value = 0.996991161888649
$ws.Range("C15").Value = 0.8095523686001933
$ws.Range("D15").Value = 0.9868251682110373
$ws.Range("E15").Value = 0.9987675056407571
$ws.Range("F15").Value = 0.9974856634705626
$ws.Range("G15").Value = 0.02012011430502972
$ws.Range("H15").Value = 1.273524188101216
$ws.Range("I15").Value = 0.02238712494987864
$ws.Range("J15").Value = 0.01707891852408543
$ws.Range("K15").Value = 0.01973303093569455
$ws.Range("L15").Value = 0.3032587262123228
$ws.Range("M15").Value = 0.1418453887337538
$ws.Range("N15").Value = 1.001604713659387
$ws.Range("O15").Value = 0.1478840405925833
$ws.Range("P15").Value = 145.8120705052045
$ws.Range("Q15").Value = 229.9145024211103
$ws.Range("B16").Value = 0.996825893804155
$ws.Range("C16").Value = 0.8091304839639802
$ws.Range("D16").Value = 0.9882855046340485
$ws.Range("E16").Value = 0.9989086973237947
$ws.Range("F16").Value = 0.9977684089959152
$ws.Range("G16").Value = 0.02122526274703057
$ws.Range("H16").Value = 1.276345332606173
$ws.Range("I16").Value = 0.01990567133479729
$ws.Range("J16").Value = 0.01512239739861856
$ws.Range("K16").Value = 0.01751398581846772
$ws.Range("L16").Value = 0.3192453326500265
$ws.Range("M16").Value = 0.1456889245860184
$ws.Range("N16").Value = 1.001692856637784
$ws.Range("O16").Value = 0.151891203723295
$ws.Range("P16").Value = 145.705126335463
$ws.Range("Q16").Value = 229.8075582513689
$ws.Range("B17").Value = 0.9966231582100686
$ws.Range("C17").Value = 0.8086582553750961
$ws.Range("D17").Value = 0.9897216569904994
$ws.Range("E17").Value = 0.9990473770749935
$ws.Range("F17").Value = 0.9980463053034152
$ws.Range("G17").Value = 0.02258095659819752
$ws.Range("H17").Value = 1.279503127354454
$ws.Range("I17").Value = 0.01746531211306788
$ws.Range("J17").Value = 0.01320068461031956
$ws.Range("K17").Value = 0.01533299836169372
$ws.Range("L17").Value = 0.3362693323426906
$ws.Range("M17").Value = 0.1502696130233838
$ws.Range("N17").Value = 1.001800982287963
$ws.Range("O17").Value = 0.1566669015507712
$ws.Range("P17").Value = 145.5812967126262
$ws.Range("Q17").Value = 229.683728628532
$ws.Range("B18").Value = 0.9963773458082854
$ws.Range("C18").Value = 0.8081298418461407
$ws.Range("D18").Value = 0.9911215523089242
$ws.Range("E18").Value = 0.9991823289193833
$ws.Range("F18").Value = 0.998316991654469
$ws.Range("G18").Value = 0.02422470526078426
$ws.Range("H18").Value = 1.2830366310557
$ws.Range("I18").Value = 0.01508656209087935
$ws.Range("J18").Value = 0.01133063016526296
$ws.Range("K18").Value = 0.01320859612807115
$ws.Range("L18").Value = 0.3544044367899156
$ws.Range("M18").Value = 0.1556428773210784
$ws.Range("N18").Value = 1.001932082235581
$ws.Range("O18").Value = 0.1622689168338091
$ws.Range("P18").Value = 145.4407645760336
$ws.Range("Q18").Value = 229.5431964919395
$ws.Range("B19").Value = 0.9960821750128371
$ws.Range("C19").Value = 0.8075387613506988
$ws.Range("D19").Value = 0.9924712553872501
$ws.Range("E19").Value = 0.999312203043291
$ws.Range("F19").Value = 0.9985777616404815
$ws.Range("G19").Value = 0.02619851372908343
$ws.Range("H19").Value = 1.286989189050395
$ws.Range("I19").Value = 0.01279310044038357
$ws.Range("J19").Value = 0.009530938698950261
$ws.Range("K19").Value = 0.01116201956966692
$ws.Range("L19").Value = 0.3737232894867413
$ws.Range("M19").Value = 0.1618595493910799
$ws.Range("N19").Value = 1.00208950665982
$ws.Range("O19").Value = 0.1687502455041159
$ws.Range("P19").Value = 145.2841051954435
$ws.Range("Q19").Value = 229.3865371113494
$ws.Range("B20").Value = 0.9957305637375972
$ws.Range("C20").Value = 0.8068777985734649
$ws.Range("D20").Value = 0.9937556727090312
$ws.Range("E20").Value = 0.9994355155263954
$ws.Range("F20").Value = 0.9988256736751288
$ws.Range("G20").Value = 0.02854973994563435
$ws.Range("H20").Value = 1.291409050185212
$ws.Range("I20").Value = 0.010610574581147
$ws.Range("J20").Value = 0.007822173189277127
$ws.Range("K20").Value = 0.009216354861801062
$ws.Range("L20").Value = 0.3943076410152645
$ws.Range("M20").Value = 0.1689666829455865
$ws.Range("N20").Value = 1.002277032673281
$ws.Range("O20").Value = 0.1761599444478325
$ws.Range("P20").Value = 145.1122149033119
$ws.Range("Q20").Value = 229.2146468192177
$ws.Range("B21").Value = 0.995314450125328
$ws.Range("C21").Value = 0.8061389615669636
$ws.Range("D21").Value = 0.9949582665947555
$ws.Range("E21").Value = 0.9995505555594983
$ws.Range("F21").Value = 0.9990574198334148
$ws.Range("G21").Value = 0.03133229358690458
$ws.Range("H21").Value = 1.296349656649717
$ws.Range("I21").Value = 0.008567085904029787
$ws.Range("J21").Value = 0.006228040658252684
$ws.Range("K21").Value = 0.007397563281141236
$ws.Range("L21").Value = 0.4162395907098929
$ws.Range("M21").Value = 0.1770093036732945
$ws.Range("N21").Value = 1.002498959933158
$ws.Range("O21").Value = 0.184544956190439
$ws.Range("P21").Value = 144.926211938655
$ws.Range("Q21").Value = 229.0286438545608
$ws.Range("B22").Value = 0.9948249059131065
$ws.Range("C22").Value = 0.8053134128705886
$ws.Range("D22").Value = 0.9960608910821209
$ws.Range("E22").Value = 0.9996555181956341
$ws.Range("F22").Value = 0.9992694511903095
$ws.Range("G22").Value = 0.03460587798817447
$ws.Range("H22").Value = 1.30187010458368
$ws.Range("I22").Value = 0.006693468649035596
$ws.Range("J22").Value = 0.004773552613587151
$ws.Range("K22").Value = 0.00573349752226014
$ws.Range("L22").Value = 0.4396093406568472
$ws.Range("M22").Value = 0.1860265518364904
$ws.Range("N22").Value = 1.002760050179677
$ws.Range("O22").Value = 0.1939460872762192
$ws.Range("P22").Value = 144.7274634545917
$ws.Range("Q22").Value = 228.8298953704975
$ws.Range("B23").Value = 0.9942517112914248
$ws.Range("C23").Value = 0.8043912496863987
$ws.Range("D23").Value = 0.9970438468973863
$ws.Range("E23").Value = 0.9997483154234428
$ws.Range("F23").Value = 0.9994577843916235
$ws.Range("G23").Value = 0.03843883306267837
$ws.Range("H23").Value = 1.308036614042527
$ws.Range("I23").Value = 0.005023196496112148
$ws.Range("J23").Value = 0.003487643042383331
$ws.Range("K23").Value = 0.00425541976924774
$ws.Range("L23").Value = 0.4645164299069868
$ws.Range("M23").Value = 0.1960582389563835
$ws.Range("N23").Value = 1.003065753977907
$ws.Range("O23").Value = 0.2044048440852612
$ws.Range("P23").Value = 144.5173741054226
$ws.Range("Q23").Value = 228.6198060213285
$ws.Range("B24").Value = 0.9935834674281447
$ws.Range("C24").Value = 0.8033616308625886
$ws.Range("D24").Value = 0.9978862661503503
$ws.Range("E24").Value = 0.9998266908255304
$ws.Range("F24").Value = 0.9996181733458339
$ws.Range("G24").Value = 0.04290738285341177
$ws.Range("H24").Value = 1.314921679858305
$ws.Range("I24").Value = 0.003591728878279841
$ws.Range("J24").Value = 0.002401579567521292
$ws.Range("K24").Value = 0.002996654222900567
$ws.Range("L24").Value = 0.491062774725785
$ws.Range("M24").Value = 0.2071409733814432
$ws.Range("N24").Value = 1.003422150704989
$ws.Range("O24").Value = 0.2159593934592188
$ws.Range("P24").Value = 144.2974227467392
$ws.Range("Q24").Value = 228.399854662645
$ws.Range("B25").Value = 0.992807329236903
$ws.Range("C25").Value = 0.80221252287599
$ws.Range("D25").Value = 0.9985657552203273
$ws.Range("E25").Value = 0.9998881154127998
$ws.Range("F25").Value = 0.9997459594769608
$ws.Range("G25").Value = 0.04809742251203331
$ws.Range("H25").Value = 1.322605770255846
$ws.Range("I25").Value = 0.002437117802000525
$ws.Range("J25").Value = 0.001550406891978315
$ws.Range("K25").Value = 0.00199376234698942
$ws.Range("L25").Value = 0.5193564912407961
$ws.Range("M25").Value = 0.2193112457491255
$ws.Range("N25").Value = 1.003836091073652
$ws.Range("O25").Value = 0.2286477795175301
$ws.Range("P25").Value = 144.0690533786692
$ws.Range("Q25").Value = 228.1714852945751
$ws.Range("B26").Value = 0.9919088644376031
$ws.Range("C26").Value = 0.8009304961334467
$ws.Range("D26").Value = 0.999058548557948
$ws.Range("E26").Value = 0.9999297897669113
$ws.Range("F26").Value = 0.9998360986176116
$ws.Range("G26").Value = 0.05410546076200164
$ws.Range("H26").Value = 1.331178688986423
$ws.Range("I26").Value = 0.001599746501896126
$ws.Range("J26").Value = 0.000972917110319929
$ws.Range("K26").Value = 0.001286331806108027
$ws.Range("L26").Value = 0.5495231894412111
$ws.Range("M26").Value = 0.232605805520846
$ws.Range("N26").Value = 1.004315272299945
$ws.Range("O26").Value = 0.2425083162222655
$ws.Range("P26").Value = 143.8336403198483
$ws.Range("Q26").Value = 227.9360722357542
